$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = 2967942.6
$ws.Range("C7").Value = -33.20075587597366
$ws.Range("D7").Value = 3022
$ws.Range("E7").Value = 3022
$ws.Range("F7").Value = 982.112045003309
$ws.Range("G7").Value = 4.686042412769265
